$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 640
$ws.Range("L3").Value = 632
$ws.Range("D4").Value = 1985
$ws.Range("I4").Value = 1822
$ws.Range("K4").Value = 1737
$ws.Range("L4").Value = 166
$ws.Range("I5").Value = 727
$ws.Range("L5").Value = 52
$ws.Range("L6").Value = 706
$ws.Range("D7").Value = 28175
$ws.Range("I7").Value = 26286
$ws.Range("K7").Value = 27528
$ws.Range("L7").Value = 2196

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L7").Value = 78
$ws.Range("L18").Value = 16
$ws.Range("L19").Value = 74
$ws.Range("L22").Value = 7
$ws.Range("L23").Value = 24
$ws.Range("L27").Value = 23
$ws.Range("L30").Value = 16
$ws.Range("L33").Value = 94
$ws.Range("L35").Value = 3
$ws.Range("L48").Value = 39
$ws.Range("L49").Value = 13
$ws.Range("L52").Value = 43
$ws.Range("L54").Value = 46
$ws.Range("D63").Value = 364
$ws.Range("I63").Value = 245
$ws.Range("K63").Value = 78
$ws.Range("L63").Value = 11
$ws.Range("L67").Value = 72
$ws.Range("L77").Value = 16
$ws.Range("L79").Value = 67
$ws.Range("L85").Value = 106
$ws.Range("L86").Value = 15
$ws.Range("L88").Value = 36
$ws.Range("L92").Value = 4
$ws.Range("L95").Value = 30
$ws.Range("L96").Value = 21
$ws.Range("D101").Value = 28175
$ws.Range("I101").Value = 26286
$ws.Range("K101").Value = 27528
$ws.Range("L101").Value = 2196

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("L2").Value = 11
$ws.Range("L4").Value = 5
$ws.Range("L7").Value = 21

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L2").Value = 17
$ws.Range("L3").Value = 27
$ws.Range("L7").Value = 78

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L3").Value = 46
$ws.Range("L7").Value = 106

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L2").Value = 15
$ws.Range("L3").Value = 12
$ws.Range("L7").Value = 43

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L2").Value = 19
$ws.Range("L3").Value = 33
$ws.Range("L7").Value = 94

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("L5").Value = 2
$ws.Range("L7").Value = 30

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("L6").Value = 6
$ws.Range("L7").Value = 16

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L2").Value = 22
$ws.Range("L4").Value = 7
$ws.Range("L7").Value = 72

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("L6").Value = 6
$ws.Range("L7").Value = 13

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("L2").Value = 15
$ws.Range("L6").Value = 22
$ws.Range("L7").Value = 46

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("L4").Value = 10
$ws.Range("L7").Value = 39

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L2").Value = 26
$ws.Range("L3").Value = 17
$ws.Range("L6").Value = 26
$ws.Range("L7").Value = 74

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("L3").Value = 12
$ws.Range("L7").Value = 24

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L3").Value = 23
$ws.Range("L5").Value = 3
$ws.Range("L6").Value = 14
$ws.Range("L7").Value = 67

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("L3").Value = 6
$ws.Range("L7").Value = 16

$ws = $wb.Worksheets.Item("Gold Coast")
$ws.Range("L4").Value = 2
$ws.Range("L7").Value = 3

$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("L6").Value = 1
$ws.Range("L7").Value = 4

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("L6").Value = 17
$ws.Range("L7").Value = 36

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("L2").Value = 7
$ws.Range("L3").Value = 7
$ws.Range("L6").Value = 6
$ws.Range("L7").Value = 23

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("L4").Value = 12
$ws.Range("L7").Value = 15

$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("L2").Value = 3
$ws.Range("L7").Value = 7

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("L6").Value = 4
$ws.Range("L7").Value = 16
